$p = $ppt.ActivePresentation

# Slide 1
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "The Assessment List for Trustworthy AI (ALTAI) is a tool designed to help foster responsible and sustainable AI innovation in Europe. It aims to make ethics a core pillar of AI development, ensuring that AI systems are lawful, ethical, and robust. ALTAI seeks to benefit, empower, and protect both individuals and society as a whole."
$tr.Font.Name = "Calibri"

# Slide 2
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI provides a structured framework to assess the trustworthiness of AI systems. It helps identify potential risks and ensure AI aligns with ethical principles and societal values. By addressing key questions about AI development and deployment, ALTAI promotes transparency and accountability in AI."
$tr.Font.Name = "Calibri"

# Slide 3
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI aims to ensure AI systems are developed and deployed responsibly. It helps identify potential risks and biases, promoting transparency and accountability in AI development."
$tr.Font.Name = "Calibri"

# Slide 4
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI is developed through a collaborative process involving experts from various fields, including AI, ethics, law, and social sciences. This ensures a comprehensive and multi-faceted approach to assessing trustworthiness in AI systems."
$tr.Font.Name = "Calibri"

# Slide 5
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI considers the impact of AI on fundamental rights. It helps ensure AI systems respect human dignity, privacy, and other essential rights."
$tr.Font.Name = "Calibri"

# Slide 6
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "AI systems should support human decision-making and uphold fundamental rights. This means AI should empower users and ensure human oversight prevents potential harm to autonomy.  ALTAI assesses AI's impact on human behavior, including its influence on decision-making processes, perception, and trust. It also considers how AI systems that mimic human behavior might affect human relationships and dependence."
$tr.Font.Name = "Calibri"

# Slide 7
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "To ensure responsible use, ALTAI emphasizes the need for human oversight mechanisms in AI systems. These mechanisms allow humans to monitor, intervene, and correct AI actions when necessary, safeguarding against unintended consequences and upholding human control."
$tr.Font.Name = "Calibri"

# Slide 8
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI recognizes the importance of preventing over-reliance on AI systems.  It promotes human oversight mechanisms to ensure humans remain in control of critical decisions.  These mechanisms help prevent AI from making decisions that could potentially harm individuals or violate their rights."
$tr.Font.Name = "Calibri"

# Slide 9
$s = $p.Slides.Item(9)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI highlights the importance of user awareness in trustworthy AI. Users should understand how AI systems work, their limitations, and potential biases."
$tr.Font.Name = "Calibri"

# Slide 10
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "AI systems should be designed to interact with humans in a way that respects social norms and avoids potential harm. This includes mitigating risks of bias, discrimination, and manipulation in social interactions."
$tr.Font.Name = "Calibri"

# Slide 11
$s = $p.Slides.Item(11)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "Could the AI system be harmed or cause harm due to risks like design flaws, attacks, or misuse?  Is the AI system secure against cyberattacks like data poisoning, model evasion, or model inversion? What steps have been taken to protect the AI system from attacks throughout its lifespan? Has the system been tested for vulnerabilities (red-teaming or penetration testing)?"
$tr.Font.Name = "Calibri"

# Slide 12
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "AI systems should be designed to be robust and safe. This includes considering potential risks like design flaws, attacks, or misuse.  It's important to ensure the AI system can handle unexpected inputs and avoid causing harm."
$tr.Font.Name = "Calibri"

# Slide 13
$s = $p.Slides.Item(13)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "The goal of an AI model is to learn patterns that generalize well for unseen data. It is important to check if a trained AI model is performing well on unseen examples that have not been used for training the model. To do this, the model is used to predict the answer on the test dataset and then the predicted target is compared to the actual answer.  The concept of accuracy is used to evaluate the predictive capability of the AI model."
$tr.Font.Name = "Calibri"

# Slide 14
$s = $p.Slides.Item(14)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "AI reliability means an AI system behaves as expected, even with new data it hasn't encountered before.  Reproducibility ensures that the same inputs consistently produce the same outputs, allowing for verification and trust in the AI system."
$tr.Font.Name = "Calibri"

# Slide 15
$s = $p.Slides.Item(15)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "How can we protect AI systems from attacks or misuse? This includes considering measures to prevent unauthorized access, data breaches, and manipulation of the AI's decision-making process."
$tr.Font.Name = "Calibri"

# Slide 16
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI asks:" + "`r" + ("" + ([char]0x2022) + " Did you consider the privacy and data protection implications of the AI system's non-personal training data or other processed non-personal data?") + "`r" + "This ensures responsible use of data throughout the AI system's lifecycle."
$tr.Font.Name = "Calibri"

# Slide 17
$s = $p.Slides.Item(17)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "This section helps assess how your AI system handles data protection." + "`r" + "* Does your AI system use personal data (including sensitive categories like health or political beliefs) during training or development?" + "`r" + "* Have you implemented measures like a Data Protection Impact Assessment (DPIA) or designated a Data Protection Officer (DPO) to ensure compliance with data protection regulations like GDPR?"
$tr.Font.Name = "Calibri"

# Slide 18
$s = $p.Slides.Item(18)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI emphasizes the importance of tracing an AI system's journey.  This includes tracking the data used by the AI, the AI model's decisions, and the system's outputs.  Effective logging practices are crucial for this traceability, allowing us to understand how the AI arrived at its conclusions."
$tr.Font.Name = "Calibri"

# Slide 19
$s = $p.Slides.Item(19)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI encourages a thorough assessment of the potential impact of your AI system on data protection. This involves considering the privacy implications throughout the AI system's lifecycle, from data collection to deployment and beyond."
$tr.Font.Name = "Calibri"

# Slide 20
$s = $p.Slides.Item(20)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI promotes the concept of 'Privacy-by-Design', ensuring privacy is integrated into every stage of an AI system's development. This means proactively considering the privacy implications of your AI system's design and functionality from the very beginning."
$tr.Font.Name = "Calibri"

# Slide 21
$s = $p.Slides.Item(21)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI stresses the importance of tracing an AI system's development and operation. This includes tracking the data used for training, the algorithms employed, and the system's decision-making processes.  Establishing clear audit trails helps ensure transparency and accountability in AI systems."
$tr.Font.Name = "Calibri"

# Slide 22
$s = $p.Slides.Item(22)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI emphasizes the importance of understanding how AI systems arrive at their decisions.  This transparency is crucial for building trust and accountability in AI systems."
$tr.Font.Name = "Calibri"

# Slide 23
$s = $p.Slides.Item(23)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI emphasizes clear communication about potential limitations of AI systems. This includes establishing mechanisms for flagging issues related to bias, discrimination, or poor performance.  It also stresses the importance of defining clear steps and communication channels for raising such issues."
$tr.Font.Name = "Calibri"

# Slide 24
$s = $p.Slides.Item(24)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI emphasizes the importance of clear communication about potential limitations of AI systems. This includes establishing transparent processes for notifying users about the system's capabilities and potential biases."
$tr.Font.Name = "Calibri"

# Slide 25
$s = $p.Slides.Item(25)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI highlights the importance of ensuring the quality of data used to train and operate AI systems. This includes assessing the accuracy, completeness, and relevance of the data."
$tr.Font.Name = "Calibri"

# Slide 26
$s = $p.Slides.Item(26)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI recognizes the potential for AI systems to perpetuate or amplify existing societal biases.  It encourages developers to consider the diverse perspectives of stakeholders throughout the AI system's lifecycle to mitigate unfair bias.  This includes actively seeking input from individuals and groups who may be disproportionately affected by the AI system."
$tr.Font.Name = "Calibri"

# Slide 27
$s = $p.Slides.Item(27)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI promotes the development of AI systems that are accessible to all individuals, regardless of their abilities or disabilities. This includes considering the needs of users with visual, auditory, or cognitive impairments."
$tr.Font.Name = "Calibri"

# Slide 28
$s = $p.Slides.Item(28)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI emphasizes the importance of involving stakeholders throughout the AI development process. This includes engaging with potential end-users and subject communities to understand their needs and perspectives. It also encourages assessing the potential impact of the AI system on different groups and identifying any risks of disproportionate or unfair effects."
$tr.Font.Name = "Calibri"

# Slide 29
$s = $p.Slides.Item(29)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI defines fairness in AI systems as the absence of discrimination and the promotion of equitable outcomes for all individuals."
$tr.Font.Name = "Calibri"

# Slide 30
$s = $p.Slides.Item(30)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI emphasizes continuous monitoring of AI systems for potential biases. This includes identifying and documenting conflicts between fairness principles and explaining any 'trade-off' decisions made.  Training is provided to those involved in this process, ensuring they understand the legal framework applicable to the AI system."
$tr.Font.Name = "Calibri"

# Slide 31
$s = $p.Slides.Item(31)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "AI systems should be developed and used in an environmentally friendly way.  This means considering the impact of the AI system's development, deployment, and use on the environment.  For example, the amount of energy used and carbon emissions should be evaluated.  Measures to reduce the environmental impact of the AI system throughout its lifecycle should be encouraged."
$tr.Font.Name = "Calibri"

# Slide 32
$s = $p.Slides.Item(32)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI recognizes the potential impact of AI on jobs and skills.  It encourages assessing how AI systems might affect employment opportunities and the skills needed in the future."
$tr.Font.Name = "Calibri"

# Slide 33
$s = $p.Slides.Item(33)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI highlights the importance of AI's impact on society and democracy.  This includes assessing how AI systems might affect  political discourse, social cohesion, and the distribution of power."
$tr.Font.Name = "Calibri"

# Slide 34
$s = $p.Slides.Item(34)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "AI systems should be developed and used in an environmentally friendly way. This means considering their energy consumption and potential impact on natural resources."
$tr.Font.Name = "Calibri"

# Slide 35
$s = $p.Slides.Item(35)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ALTAI emphasizes the importance of involving workers in the development and deployment of AI systems. This ensures that their perspectives and concerns are considered, promoting fairness and transparency."
$tr.Font.Name = "Calibri"

# Slide 36
$s = $p.Slides.Item(36)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "Trustworthy AI requires mechanisms to allow for independent audits. These audits should enable the evaluation of AI systems' compliance with ethical principles and legal requirements."
$tr.Font.Name = "Calibri"

# Slide 37
$s = $p.Slides.Item(37)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "Trustworthy AI demands robust risk management processes. These processes should identify, assess, and mitigate potential risks associated with AI systems throughout their lifecycle."
$tr.Font.Name = "Calibri"

# Slide 38
$s = $p.Slides.Item(38)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "Trustworthy AI necessitates independent audits conducted by third parties. These audits ensure transparency and allow for the evaluation of AI systems' adherence to ethical and safety standards."
$tr.Font.Name = "Calibri"

# Slide 39
$s = $p.Slides.Item(39)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "Trustworthy AI necessitates the establishment of ethics review boards. These boards should evaluate AI systems for potential biases, fairness, and societal impacts."
$tr.Font.Name = "Calibri"

# Slide 40
$s = $p.Slides.Item(40)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "Trustworthy AI should incorporate mechanisms for redress. This means designing systems that allow for the fair and effective resolution of complaints or harms caused by AI."
$tr.Font.Name = "Calibri"

